$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description text in B2: remove the "/RME" segment
$old = $ws.Range("B2").Value2
$new = $old -replace "4% S/LFM\+CDL/RME/H:1", "4% S/LFM+CDL/H:1"
$ws.Range("B2").Value = $new

# Apply wrap text to B2 and resize row 2 to the tallest row height
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 409.6

# Update the active selection to A7, matching the saved view state
$ws.Range("A7").Select()

$wb.Save()
